$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.416.38'
$ws.Range("E2").Value = '  -1.97%  '
$ws.Range("D3").Value = '1.839.81'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.45'
$ws.Range("E5").Value = '  -6.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5224'
$ws.Range("E7").Value = '  -1.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3240'
$ws.Range("E8").Value = '  -6.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06795'
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("E10").Value = '  -7.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7650'
$ws.Range("E11").Value = '  -5.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07681'
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("D13").Value = '1.841.32'
$ws.Range("E13").Value = '  -1.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.49'
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.023'
$ws.Range("E15").Value = '  -3.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.94'
$ws.Range("E17").Value = '  -4.78%  '
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007911'
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("D20").Value = '26.457.55'
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("D21").Value = '2.075.35'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.571'
$ws.Range("E22").Value = '  -3.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.452'
$ws.Range("E23").Value = '  -6.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.944'
$ws.Range("E24").Value = '  -4.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.77'
$ws.Range("E25").Value = '  -1.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.211'
$ws.Range("E26").Value = '  -6.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.651'
$ws.Range("E27").Value = '  -0.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.95'
$ws.Range("E28").Value = '  -2.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.46'
$ws.Range("E29").Value = '  -2.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.182'
$ws.Range("E30").Value = '  -4.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.141'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08740'
$ws.Range("E32").Value = '  -1.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04813'
$ws.Range("E33").Value = '  -2.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.124'
$ws.Range("E34").Value = '  -5.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.846'
$ws.Range("E35").Value = '  -1.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7027'
$ws.Range("E36").Value = '  -5.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.067'
$ws.Range("E37").Value = '  -6.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01766'
$ws.Range("E38").Value = '  -4.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.186'
$ws.Range("E39").Value = '  -9.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4826'
$ws.Range("E40").Value = '  -6.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '111.45'
$ws.Range("E41").Value = '  -4.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8882'
$ws.Range("E42").Value = '  -7.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.096'
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.0000'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.661'
$ws.Range("E45").Value = '  -5.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4132'
$ws.Range("E46").Value = '  -8.84%  '
$ws.Range("E47").Value = '  -1.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.009'
$ws.Range("E48").Value = '  -3.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.78'
$ws.Range("E49").Value = '  -4.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1221'
$ws.Range("E50").Value = '  -9.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8821'
$ws.Range("E51").Value = '  -0.62%  '
